$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume 1h (E) columns to match refreshed symbol data.
# Force text storage (matching original inlineStr text cells) via Text number
# format, then strip the resulting style delta back to Normal so no stray
# cell-style attribute is introduced.
$cellUpdates = @{
    "D2" = "309.10"
    "E2" = "2.16%"
    "D3" = "38.77"
    "E3" = "8.68%"
    "E4" = "1.14%"
    "D5" = "0.08174"
    "E5" = "3.25%"
    "D6" = "1.987"
    "E6" = "6.65%"
    "D7" = "7.898"
    "E7" = "1.50%"
    "D8" = "0.9319"
    "E8" = "1.46%"
    "E9" = "4.69%"
    "D10" = "0.1951"
    "E10" = "3.34%"
    "D11" = "0.09187"
    "E11" = "1.37%"
    "D12" = "0.03445"
    "E12" = "-0.07%"
    "D13" = "0.09859"
    "E13" = "0.58%"
    "D14" = "0.001411"
    "E14" = "-0.11%"
    "D15" = "0.006211"
    "E15" = "3.23%"
    "D16" = "3.760"
    "E16" = "0.85%"
    "D17" = "4.185"
    "E17" = "1.96%"
    "D18" = "3.489"
    "E18" = "5.73%"
    "E19" = "0.35%"
    "D20" = "0.1305"
    "E20" = "0.32%"
    "D21" = "4.814"
    "E21" = "-6.75%"
    "D22" = "0.2455"
    "E22" = "11.97%"
    "D23" = "0.04475"
    "E23" = "1.64%"
    "D24" = "0.001238"
    "E24" = "0.20%"
    "E25" = "-9.35%"
    "E27" = "0.15%"
    "D39" = "0.02119"
    "E39" = "10.13%"
    "D40" = "0.05174"
    "E40" = "-1.42%"
    "D41" = "0.007467"
    "E41" = "-2.08%"
    "D42" = "0.009963"
    "E42" = "-1.98%"
    "D43" = "0.1368"
    "E43" = "1.99%"
    "E44" = "-0.32%"
    "D45" = "0.009688"
    "E45" = "-4.37%"
    "D46" = "0.00006313"
    "E46" = "2.62%"
    "E47" = "0.20%"
    "E48" = "-0.63%"
    "D49" = "0.001603"
    "E49" = "-3.35%"
    "D50" = "0.00002104"
    "E50" = "0.20%"
    "D51" = "0.0002004"
    "E51" = "0.20%"
}

foreach ($addr in $cellUpdates.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $cellUpdates[$addr]
    $rng.Style = "Normal"
}
